$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update meanEMG / legmaxROM values (row 1 headers + row 2/3 data) for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 285.29093073464242
$ws.Range("C2").Value = 264.61263669867992
$ws.Range("D2").Value = 284.30665607591612
$ws.Range("E2").Value = 262.21989162037539

$ws.Range("B3").Value = 308.9846248534343
$ws.Range("C3").Value = 261.16854873030132
$ws.Range("D3").Value = 315.01383289457624
$ws.Range("E3").Value = 255.34581346004452

# Update the active selection to match the edited range
$ws.Range("B1:E3").Select() | Out-Null
